# Update cryptos list data (prices and 1h volume changes), and for row 51
# replace RenderToken with Aptos (coin name, link, price, volume change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{Row=2; D='28.967.23'; E='  -0.27%  '},
  @{Row=3; D='1.830.13'; E='  -0.06%  '},
  @{Row=4; D='0.9963'; E='  -0.28%  '},
  @{Row=5; D='241.63'; E='  -0.02%  '},
  @{Row=6; D='0.6267'; E='  -4.23%  '},
  @{Row=7; D='0.9972'; E='  -0.31%  '},
  @{Row=8; D='0.07591'; E='  +3.52%  '},
  @{Row=9; D='0.2922'; E='  -0.44%  '},
  @{Row=10; D='22.53'; E='  -1.73%  '},
  @{Row=11; D='0.07712'; E='  +0.64%  '},
  @{Row=12; D='1.834.25'; E='  +0.04%  '},
  @{Row=13; D='4.942'; E='  -0.73%  '},
  @{Row=14; D='0.6640'; E='  -0.38%  '},
  @{Row=15; D='0.00001027'; E='  +17.84%  '},
  @{Row=16; D='82.72'; E='  +1.18%  '},
  @{Row=17; D='6.041'; E='  -1.06%  '},
  @{Row=18; D='28.998.91'; E='  -0.15%  '},
  @{Row=19; D='226.26'; E='  +1.10%  '},
  @{Row=20; D='12.32'; E='  -0.87%  '},
  @{Row=21; D='0.9961'; E='  -0.43%  '},
  @{Row=22; D='7.176'; E='  +0.87%  '},
  @{Row=23; D='0.9971'; E='  -0.26%  '},
  @{Row=24; D='158.14'; E='  +0.37%  '},
  @{Row=25; D='8.485'; E='  -0.04%  '},
  @{Row=26; D='0.1373'; E='  -0.34%  '},
  @{Row=27; D='17.88'; E='  -0.17%  '},
  @{Row=28; D='1.486'; E='  -1.02%  '},
  @{Row=29; D='4.097'; E='  -0.16%  '},
  @{Row=30; D='4.010'; E='  -0.05%  '},
  @{Row=31; D='1.186'; E='  -1.10%  '},
  @{Row=32; D='0.05221'; E='  -2.37%  '},
  @{Row=33; D='1.839'; E='  +0.06%  '},
  @{Row=34; D='0.7346'; E='  -1.04%  '},
  @{Row=35; D='1.136'; E='  -1.72%  '},
  @{Row=36; D='2.689'; E='  +1.51%  '},
  @{Row=37; D='1.234.55'; E='  -4.90%  '},
  @{Row=38; D='2.752'; E='  +0.15%  '},
  @{Row=39; E='  -0.15%  '},
  @{Row=40; D='6.354'; E='  +0.48%  '},
  @{Row=41; D='0.8940'; E='  -0.53%  '},
  @{Row=42; D='0.9977'; E='  -0.23%  '},
  @{Row=43; D='101.65'; E='  -1.68%  '},
  @{Row=44; D='1.982.64'; E='  -0.21%  '},
  @{Row=45; D='0.00000000124'; E='  +3.32%  '},
  @{Row=46; D='63.99'; E='  -0.54%  '},
  @{Row=47; D='0.5090'; E='  -0.85%  '},
  @{Row=48; D='0.4034'; E='  +1.35%  '},
  @{Row=49; D='8.907'; E='  +2.70%  '},
  @{Row=50; D='0.05749'; E='  -1.47%  '},
  @{Row=51; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='6.685'; E='  -0.36%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey('D')) {
        # Force text storage so numeric-looking price strings (which may
        # contain multiple dots, or significant trailing zeros) keep their
        # exact original textual representation instead of being coerced
        # into a floating point number by Excel's automatic type detection.
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey('E')) { $ws.Range("E$r").Value = $u.E }
}
